# Commit: "feat: add 2022-Q4 data"
#
# 1) Insert a new worksheet "2022-Q4" right after "总计" (i.e. before the
#    sheet that is currently "2022-Q3"), populate it with the fund detail
#    table for the new quarter, modelled on the layout of the existing
#    quarter sheets.
# 2) On the "总计" summary sheet, insert the new 2022-Q4 summary row at the
#    top of the data (row 2) and push the existing rows down by one.

$wb = $excel.ActiveWorkbook

# --- Grab stable references before the sheet collection gets reshuffled ---
$total = $wb.Worksheets.Item(1)      # "总计" - always first, index is stable
$q3anchor = $wb.Worksheets.Item(2)   # currently "2022-Q3" - insertion point

# --- 1) Create the new "2022-Q4" detail sheet, positioned before 2022-Q3 ---
$newSheet = $wb.Worksheets.Add($q3anchor)
$newSheet.Name = "2022-Q4"

# Clone the layout/styling of the existing "2022-Q3" detail sheet (header
# row + indexed first column formatting) by copying its first 4 rows, then
# overwrite with the 2022-Q4 values.
$q3 = $wb.Worksheets.Item("2022-Q3")
$q3.Range("A1:H4").Copy($newSheet.Range("A1:H4"))

# Header row (unchanged wording, just re-asserting values over the copy)
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# Data rows (fund code / D-G numeric-looking figures kept as text, like the
# source data; H is a genuine number)
$newSheet.Range("A2").Value = 0
$newSheet.Range("B2").Value = "'952004"
$newSheet.Range("C2").Value = "国泰君安君得明混合"
$newSheet.Range("D2").Value = "'18.25"
$newSheet.Range("E2").Value = "'83.43"
$newSheet.Range("F2").Value = "'2.29"
$newSheet.Range("G2").Value = "'0.4179"
$newSheet.Range("H2").Value = 5

$newSheet.Range("A3").Value = 1
$newSheet.Range("B3").Value = "'001899"
$newSheet.Range("C3").Value = "东海中证社会发展安全产业主题指数"
$newSheet.Range("D3").Value = "'0.20"
$newSheet.Range("E3").Value = "'93.35"
$newSheet.Range("F3").Value = "'4.04"
$newSheet.Range("G3").Value = "'0.0081"
$newSheet.Range("H3").Value = 3

$newSheet.Range("A4").Value = 2
$newSheet.Range("B4").Value = "'005616"
$newSheet.Range("C4").Value = "东方量化成长灵活配置混合"
$newSheet.Range("D4").Value = "'0.16"
$newSheet.Range("E4").Value = "'93.18"
$newSheet.Range("F4").Value = "'2.00"
$newSheet.Range("G4").Value = "'0.0032"
$newSheet.Range("H4").Value = 6

# --- 2) Update the "总计" summary sheet ---
# The sheet already has rows 2-9 holding 2022-Q3 .. 2020-Q4. We need to push
# all of that down by one row and put the new 2022-Q4 figures in row 2.
# Row-Insert() would be the obvious move, but it drags the header row's
# bold/border styling onto the freshly inserted cells (Excel's usual
# "inherit format from above" behaviour), which does not match the source
# data (plain, unstyled B/C/D cells). Instead, write every row's final
# target content directly - row 10 is new, so its index cell (column A)
# needs the same per-row style as the others, cloned from A9 first.
$total.Range("A9").Copy($total.Range("A10"))

$total.Range("B10").Value = "2020-Q4"
$total.Range("C10").Value = 7
$total.Range("D10").Value = 5.06

$total.Range("B9").Value = "2021-Q1"
$total.Range("C9").Value = 6
$total.Range("D9").Value = 14.83

$total.Range("B8").Value = "2021-Q2"
$total.Range("C8").Value = 13
$total.Range("D8").Value = 18.83

$total.Range("B7").Value = "2021-Q3"
$total.Range("C7").Value = 5
$total.Range("D7").Value = 3.05

$total.Range("B6").Value = "2021-Q4"
$total.Range("C6").Value = 3
$total.Range("D6").Value = 2.71

$total.Range("B5").Value = "2022-Q1"
$total.Range("C5").Value = 2
$total.Range("D5").Value = 0.91

$total.Range("B4").Value = "2022-Q2"
$total.Range("C4").Value = 2
$total.Range("D4").Value = 0.5

$total.Range("B3").Value = "2022-Q3"
$total.Range("C3").Value = 10
$total.Range("D3").Value = 1.1

$total.Range("B2").Value = "2022-Q4"
$total.Range("C2").Value = 3
$total.Range("D2").Value = 0.43

# Re-sequence the 0-based index in column A for every row now that row 2 is
# the new 2022-Q4 entry.
$total.Range("A2").Value = 0
$total.Range("A3").Value = 1
$total.Range("A4").Value = 2
$total.Range("A5").Value = 3
$total.Range("A6").Value = 4
$total.Range("A7").Value = 5
$total.Range("A8").Value = 6
$total.Range("A9").Value = 7
$total.Range("A10").Value = 8
